$wb = $excel.ActiveWorkbook

# Rename "wt" -> "wt_log2_expression" and "dcin5" -> "dcin5_log2_expression"
$wsWt = $wb.Worksheets.Item("wt")
$wsWt.Name = "wt_log2_expression"

$wsDcin5 = $wb.Worksheets.Item("dcin5")
$wsDcin5.Name = "dcin5_log2_expression"

# Make "wt_log2_expression" the active sheet/tab, with the given selection.
# (This also clears tabSelected on whichever sheet was previously active,
# i.e. "optimization_parameters".)
$wsWt.Activate()
$wsWt.Range("E34").Select()
